$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing prediction values
$ws.Range("B2").Value = 0.05261340675226717
$ws.Range("C2").Value = 0.9984527294684541
$ws.Range("D2").Value = 0.1677042766834937

# Add new "Modelo" column header, matching style of other headers
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$modelText = "Pipeline(steps=[('model'," + [char]10 + "                 RandomForestRegressor(max_depth=3, n_estimators=150))])"
$ws.Range("F2").Value = $modelText
